# Add two new columns (I: "I0", J: "IF") with header styling matching
# the existing header row, and fill in the corresponding data values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the new columns, copying the style of the existing header H1
# so that I1/J1 match the bold/centered/bordered look of the other headers.
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for the new I and J columns, row by row.
$data = @{
    2  = @(4, 4)
    3  = @(7, 7)
    4  = @(8, 8)
    5  = @(6, 6)
    6  = @(3, 4)
    7  = @(4, 5)
    8  = @(7, 8)
    9  = @(5, 7)
    10 = @(4, 6)
    11 = @(4, 5)
    12 = @(6, 7)
    13 = @(4, 6)
    14 = @(9, 9)
    15 = @(2, 2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]  # Column I
    $ws.Cells.Item($row, 10).Value = $vals[1] # Column J
}
